$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking prices/volumes are stored as plain text in this sheet (note
# values like "27.821.64" that are not valid numbers). For new values that
# DO look like ordinary numbers (e.g. "1.005"), Excel would normally auto-
# convert the typed text into a real number, so we briefly mark those cells
# as Text before writing them, then clear the temporary formatting again
# (the value already stored as text is unaffected by that).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "27.821.64"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "1.768.09"
$ws.Range("E3").Value = "  -2.54%  "

Set-TextValue $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.80%  "

Set-TextValue $ws.Range("D5") "338.32"
$ws.Range("E5").Value = "  +0.17%  "

Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.61%  "

Set-TextValue $ws.Range("D7") "0.3782"
$ws.Range("E7").Value = "  -3.81%  "

Set-TextValue $ws.Range("D8") "0.3380"
$ws.Range("E8").Value = "  -3.15%  "

Set-TextValue $ws.Range("D9") "45.70"
$ws.Range("E9").Value = "  -5.22%  "

Set-TextValue $ws.Range("D10") "1.132"
$ws.Range("E10").Value = "  -5.90%  "

Set-TextValue $ws.Range("D11") "0.07293"
$ws.Range("E11").Value = "  -4.09%  "

Set-TextValue $ws.Range("D12") "23.28"
$ws.Range("E12").Value = "  +4.94%  "

Set-TextValue $ws.Range("D13") "1.003"
$ws.Range("E13").Value = "  +0.87%  "

Set-TextValue $ws.Range("D14") "6.280"

Set-TextValue $ws.Range("D15") "7.294"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").Value = "1.767.99"
$ws.Range("E16").Value = "  -2.44%  "

Set-TextValue $ws.Range("D17") "0.00001061"
$ws.Range("E17").Value = "  -4.32%  "

Set-TextValue $ws.Range("D18") "0.06627"
$ws.Range("E18").Value = "  -0.85%  "

Set-TextValue $ws.Range("D19") "81.43"
$ws.Range("E19").Value = "  -4.62%  "

Set-TextValue $ws.Range("D20") "1.003"
$ws.Range("E20").Value = "  +0.68%  "

Set-TextValue $ws.Range("D21") "17.22"
$ws.Range("E21").Value = "  -3.95%  "

Set-TextValue $ws.Range("D22") "6.351"
$ws.Range("E22").Value = "  -3.53%  "

$ws.Range("D23").Value = "27.844.59"
$ws.Range("E23").Value = "  -0.20%  "

Set-TextValue $ws.Range("D24") "11.87"
$ws.Range("E24").Value = "  -7.55%  "

Set-TextValue $ws.Range("D25") "2.382"
$ws.Range("E25").Value = "  -0.82%  "

Set-TextValue $ws.Range("D26") "1.486"
$ws.Range("E26").Value = "  +0.07%  "

Set-TextValue $ws.Range("D27") "20.23"
$ws.Range("E27").Value = "  -5.41%  "

Set-TextValue $ws.Range("D28") "151.28"
$ws.Range("E28").Value = "  -2.29%  "

Set-TextValue $ws.Range("D29") "2.354"
$ws.Range("E29").Value = "  -8.30%  "

$ws.Range("D30").Value = "1.969.37"
$ws.Range("E30").Value = "  -2.59%  "

Set-TextValue $ws.Range("D31") "133.06"
$ws.Range("E31").Value = "  -1.95%  "

Set-TextValue $ws.Range("D32") "4.034"
$ws.Range("E32").Value = "  +0.13%  "

Set-TextValue $ws.Range("D33") "5.934"
$ws.Range("E33").Value = "  -3.10%  "

Set-TextValue $ws.Range("D34") "0.08761"
$ws.Range("E34").Value = "  -0.91%  "

Set-TextValue $ws.Range("D35") "12.42"
$ws.Range("E35").Value = "  -6.69%  "

Set-TextValue $ws.Range("D36") "0.02361"
$ws.Range("E36").Value = "  -3.04%  "

Set-TextValue $ws.Range("D37") "0.6725"
$ws.Range("E37").Value = "  -3.00%  "

Set-TextValue $ws.Range("D38") "0.06263"
$ws.Range("E38").Value = "  -3.93%  "

Set-TextValue $ws.Range("D39") "5.193"
$ws.Range("E39").Value = "  -6.38%  "

Set-TextValue $ws.Range("D40") "0.2129"
$ws.Range("E40").Value = "  -4.83%  "

Set-TextValue $ws.Range("D41") "1.483"
$ws.Range("E41").Value = "  -7.79%  "

Set-TextValue $ws.Range("D42") "1.217"
$ws.Range("E42").Value = "  -4.08%  "

Set-TextValue $ws.Range("D43") "8.130"
$ws.Range("E43").Value = "  -5.70%  "

Set-TextValue $ws.Range("D44") "1.002"
$ws.Range("E44").Value = "  +0.69%  "

Set-TextValue $ws.Range("D45") "14.00"
$ws.Range("E45").Value = "  -4.92%  "

Set-TextValue $ws.Range("D46") "0.6139"
$ws.Range("E46").Value = "  -6.71%  "

Set-TextValue $ws.Range("D47") "3.836"
$ws.Range("E47").Value = "  -0.88%  "

Set-TextValue $ws.Range("D48") "131.73"
$ws.Range("E48").Value = "  -0.73%  "

Set-TextValue $ws.Range("D49") "2.036"
$ws.Range("E49").Value = "  -5.92%  "

Set-TextValue $ws.Range("D50") "0.07292"
$ws.Range("E50").Value = "  +1.23%  "

Set-TextValue $ws.Range("D51") "1.187"
$ws.Range("E51").Value = "  +2.36%  "
